$d = $word.ActiveDocument

# 1) First title: "Sprint 01" -> "Initial Sprint"
$d.Content.Find.Execute("Sprint 01", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Initial Sprint", 2)

# 2) Second title: "Sprint 02" -> "Sprint 01"
$d.Content.Find.Execute("Sprint 02", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sprint 01", 2)

# 3) The single "TODO" bullet becomes the first of three planning bullets,
#    and two further bullets get appended after it in the same list.
$todoPara = $d.Paragraphs.Last
$todoPara.Range.Text = "Alain finalisiert das Design"

# The _GoBack bookmark used to sit right after "TODO"; drop it here, it is
# recreated after the last of the new bullets below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p2 = $todoPara.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Pascal & Seraphin erstellen UseCases"

$p3 = $p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
# Append a trailing placeholder character so the bookmark range below does
# not land exactly on the paragraph mark; it is stripped again immediately
# after, leaving the bookmark collapsed right after the visible text.
$p3.Range.Text = "Alle bestehen den Basic Test B in diesem SprintX"

$endPos = $p3.Range.End - 1
$placeholder = $d.Range($endPos - 1, $endPos)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder = $d.Range($endPos - 1, $endPos)
$placeholder.Delete()
